$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 357-358; existing rows 357-377 shift down to 359-379.
$ws.Range("A357:A358").EntireRow.Insert()

# --- New row 357 ---
$ws.Range("A357").Value = 10
$ws.Range("B357").Value = "Vega Modelo de Temuco"
$ws.Range("C357").Value = "La Araucanía"
$ws.Range("D357").Value = 45267
$ws.Range("E357").Value = 9
$ws.Range("F357").Value = "Fruta"
$ws.Range("G357").Value = 100101
$ws.Range("H357").Value = "Berries"
$ws.Range("I357").Value = 100112025
$ws.Range("J357").Value = "Frutilla"
$ws.Range("K357").Value = "Sin especificar"
$ws.Range("L357").Value = "Primera"
$ws.Range("M357").Value = 380
$ws.Range("N357").Value = 10000
$ws.Range("O357").Value = 10000
$ws.Range("P357").Value = 10000
$ws.Range("Q357").Value = "`$/bandeja 7 kilos"
$ws.Range("R357").Value = "Provincia de Melipilla"
$ws.Range("S357").Value = 1429
$ws.Range("T357").Value = 7

# --- New row 358 ---
$ws.Range("A358").Value = 10
$ws.Range("B358").Value = "Vega Modelo de Temuco"
$ws.Range("C358").Value = "La Araucanía"
$ws.Range("D358").Value = 45267
$ws.Range("E358").Value = 9
$ws.Range("F358").Value = "Fruta"
$ws.Range("G358").Value = 100101
$ws.Range("H358").Value = "Berries"
$ws.Range("I358").Value = 100112025
$ws.Range("J358").Value = "Frutilla"
$ws.Range("K358").Value = "Sin especificar"
$ws.Range("L358").Value = "Primera"
$ws.Range("M358").Value = 110
$ws.Range("N358").Value = 10000
$ws.Range("O358").Value = 10000
$ws.Range("P358").Value = 10000
$ws.Range("Q358").Value = "`$/caja 7 kilos"
$ws.Range("R358").Value = "Provincia de Cautín"
$ws.Range("S358").Value = 1429
$ws.Range("T358").Value = 7
